$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Review")

# Rows 3 and 4 of the review table (spreadsheet rows 4 and 5) were filled
# in with the same review metadata (No., date, status, error type, part)
# as the rows above them.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 44547
$ws.Range("C4").Value = "Open"
$ws.Range("D4").Value = "Others"
$ws.Range("E4").Value = "Improvement"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 44547
$ws.Range("C5").Value = "Open"
$ws.Range("D5").Value = "Others"
$ws.Range("E5").Value = "Improvement"

# The user scrolled back to column A and selected the merged "Error"
# cell on row 5 instead of the "Review Content" cell.
$ws.Activate() | Out-Null
$ws.Range("F5:H5").Select() | Out-Null
